$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($r1, $r2) {
    # Swap the match-detail columns (F..V) between two data rows, leaving
    # the leading Indice/pais/torneio/temporada/data_partida columns (A..E)
    # untouched, since those stay keyed to the row number.
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $addr1 = $col + $r1
        $addr2 = $col + $r2
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

Swap-RowData 74 75
Swap-RowData 84 85
Swap-RowData 106 107

# Append the new match row (row 109) after the previous last row (108),
# carrying over the same cell formatting used by the rest of the table.
$ws.Range("A108").Copy() | Out-Null
$ws.Range("A109").PasteSpecial(-4122) | Out-Null
$ws.Range("E108").Copy() | Out-Null
$ws.Range("E109").PasteSpecial(-4122) | Out-Null

$ws.Range("A109").Value2 = 108
$ws.Range("B109").Value2 = "turkey"
$ws.Range("C109").Value2 = "1-lig"
$ws.Range("D109").Value2 = "2023-2024"
$ws.Range("E109").Value2 = 45242.70833333334
$ws.Range("F109").Value2 = "Kocaelispor"
$ws.Range("G109").Value2 = 2
$ws.Range("H109").Value2 = "Corum"
$ws.Range("I109").Value2 = 1
$ws.Range("J109").Value2 = 1.93
$ws.Range("K109").Value2 = "05/11/2023 17:12"
$ws.Range("L109").Value2 = 2
$ws.Range("M109").Value2 = "12/11/2023 16:06"
$ws.Range("N109").Value2 = 3.44
$ws.Range("O109").Value2 = "05/11/2023 17:12"
$ws.Range("P109").Value2 = 3.47
$ws.Range("Q109").Value2 = "12/11/2023 16:06"
$ws.Range("R109").Value2 = 4.06
$ws.Range("S109").Value2 = "05/11/2023 17:12"
$ws.Range("T109").Value2 = 3.89
$ws.Range("U109").Value2 = "12/11/2023 16:06"
$ws.Range("V109").Value2 = "https://www.betexplorer.com/football/turkey/1-lig/kocaelispor-corum-fk/p4EkAAAi/"
